$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: K2 "Result" -> "Shown" ---
$ws.Range("K2").Value = "Shown"

# --- Existing rows 3 & 4: bump the dates forward by 3 days (keep everything else) ---
$ws.Range("F3").Value = 43817.4326388889
$ws.Range("G3").Value = 43817.4340277778
$ws.Range("H3").Value = 43817.5652777778
$ws.Range("I3").Value = 43817.5298611111
$ws.Range("J3").Value = 43817.0847222222

$ws.Range("F4").Value = 43817.4743055556
$ws.Range("G4").Value = 43817.4756944444
$ws.Range("H4").Value = 43817.6069444444
$ws.Range("I4").Value = 43817.5715277778
$ws.Range("J4").Value = 43817.0840277778

# --- New row 5 ---
$ws.Range("B5").Value = "ABC"
$ws.Range("C5").Value = "Repaired"
$ws.Range("F5").Value = 43817.2791666667
$ws.Range("F5").NumberFormat = "m/d/yy h:mm"
$ws.Range("G5").Value = 43817.28125
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = 43817.2833333333
$ws.Range("H5").NumberFormat = "m/d/yy h:mm"
$ws.Range("I5").Value = 43817.0020833333
$ws.Range("I5").NumberFormat = "m/d/yy h:mm"
$ws.Range("J5").Value = 43817.0041666667
$ws.Range("J5").NumberFormat = "m/d/yy h:mm"
$ws.Range("K5").Value = "NG"
$ws.Range("L5").Value = "F"

# --- New row 6 ---
$ws.Range("B6").Value = "DSFDSF"
$ws.Range("C6").Value = "Buy off sample"
$ws.Range("F6").Value = 43817.28125
$ws.Range("F6").NumberFormat = "m/d/yy h:mm"
